$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("GlobalConstantIntTable")

# --- Column A (id|String) for the new rows, filled top-to-bottom first ---
# (matches the shared-string insertion order seen in the target workbook)
$ws.Range("A47").Value = "Ev13CountLimit"
$ws.Range("A48").Value = "Ev14CountLimit"
$ws.Range("A49").Value = "Ev15CountLimit"
$ws.Range("A50").Value = "Ev16CountLimit"
$ws.Range("A51").Value = "Ev17CountLimit"

# --- Column C (비고/remark) for the new rows, filled top-to-bottom second ---
$ws.Range("C47").Value = "보유 스펠 종류 제한"
$ws.Range("C48").Value = "미보유 스펠 종류 제한"
$ws.Range("C49").Value = "보유 동료 종류 제한"
$ws.Range("C50").Value = "보유 동료 피피 종류 제한"
$ws.Range("C51").Value = "미보유 동료 종류 제한"

# --- Column D (value|Int) for the new rows ---
$ws.Range("D47").Value = 4
$ws.Range("D48").Value = 8
$ws.Range("D49").Value = 9
$ws.Range("D50").Value = 12
$ws.Range("D51").Value = 17

# --- Formatting: row 47 uses the same "section header" style as A32/A33,
#     rows 48-51 use the same style as A22/A23 ---
$ws.Cells.Item(32,1).Copy()
$ws.Range("A47").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Cells.Item(22,1).Copy()
$ws.Range("A48:A51").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- View state: scroll down so the newly-added rows are visible, and
#     move the active selection just past the last row of data ---
$ws.Range("A52").Select() | Out-Null
try { $excel.ActiveWindow.ScrollRow = 36 } catch {}
